$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 305, shifting existing rows 305-358 down to 306-359
$ws.Rows("305:305").Insert()

# Fill in the values for the newly inserted row 305
$ws.Cells.Item(305, 1).Value = 11
$ws.Cells.Item(305, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(305, 3).Value = "Bíobío"
$ws.Cells.Item(305, 4).Value = 44637
$ws.Cells.Item(305, 5).Value = 8
$ws.Cells.Item(305, 6).Value = "Fruta"
$ws.Cells.Item(305, 7).Value = 100104
$ws.Cells.Item(305, 8).Value = "Frutos de pepita"
$ws.Cells.Item(305, 9).Value = 100104005
$ws.Cells.Item(305, 10).Value = "Pera"
$ws.Cells.Item(305, 11).Value = "Packham's Triumph"
$ws.Cells.Item(305, 12).Value = "Primera"
$ws.Cells.Item(305, 13).Value = 220
$ws.Cells.Item(305, 14).Value = 9500
$ws.Cells.Item(305, 15).Value = 10000
$ws.Cells.Item(305, 16).Value = 9727
$ws.Cells.Item(305, 17).Value = "`$/caja 16 kilos empedrada"
$ws.Cells.Item(305, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(305, 19).Value = 608
$ws.Cells.Item(305, 20).Value = 16
